# Populate previously-empty memory-usage cells with computed values
# (changed implementation for memory calculation).
# Setting .Value on these cells also implicitly clears their placeholder
# fill/border style back to "Normal", matching the target workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.12734222412109375
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = 0.11981964111328125
$ws.Range("D2").Style = "Normal"
$ws.Range("G2").Value = 0.0
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = 1821456.0
$ws.Range("H2").Style = "Normal"
$ws.Range("A3").Value = -0.34369659423828125
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = 1.2271499633789062
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = 0.0
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = 1.2700424194335938
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = 41018832.0
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = -20456560.0
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = 0.0
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = 1861304.0
$ws.Range("H3").Style = "Normal"
$ws.Range("I3").Value = 0.0
$ws.Range("I3").Style = "Normal"
$ws.Range("J3").Value = 1331736.0
$ws.Range("J3").Style = "Normal"
$ws.Range("A4").Value = 2.91552734375
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = 2.6513214111328125
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = 1.7562789916992188
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = 1.7562789916992188
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = 6.134796142578125
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = 6.1512298583984375
$ws.Range("F4").Style = "Normal"
$ws.Range("I4").Value = 0.0
$ws.Range("I4").Style = "Normal"
$ws.Range("J4").Value = 1331736.0
$ws.Range("J4").Style = "Normal"
$ws.Range("A5").Value = 1.905059814453125
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = 3.1751022338867188
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = 1.224334716796875
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = 1.735260009765625
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = 9.676353454589844
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = 9.676322937011719
$ws.Range("F5").Style = "Normal"
$ws.Range("I5").Value = 0.0
$ws.Range("I5").Style = "Normal"
$ws.Range("J5").Value = 0.0
$ws.Range("J5").Style = "Normal"
$ws.Range("K5").Value = 0.0
$ws.Range("K5").Style = "Normal"
$ws.Range("L5").Value = 0.0
$ws.Range("L5").Style = "Normal"
$ws.Range("A6").Value = 3266064.0
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = 2770920.0
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = 1.735260009765625
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = 1.735260009765625
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = 9.576286315917969
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = 9.5762939453125
$ws.Range("F6").Style = "Normal"
$ws.Range("C7").Value = 1.735260009765625
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = 1.735260009765625
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = 10100424.0
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = 10100392.0
$ws.Range("F7").Style = "Normal"
$ws.Range("C8").Value = 665888.0
$ws.Range("C8").Style = "Normal"
$ws.Range("C9").Value = 0.0
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = 1331736.0
$ws.Range("D9").Style = "Normal"
$ws.Range("C10").Value = 1832400.0
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = 1832400.0
$ws.Range("D10").Style = "Normal"
